# Update handback status timestamps to reflect the regenerated report.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 07b0646c-...
$overview.Range("G3").Value = "2016-08-22 04:55:07"

# zh-cn sheet: Correspond Handoff / Handback datetimes for 07b0646c-...
$zhcn.Range("H3").Value = "2016-08-22 04:54:58"
$zhcn.Range("K3").Value = "2016-08-22 04:55:27"

# de-de sheet: Latest HO Xliff Generate Date + Correspond Handback datetime for 07b0646c-...
$dede.Range("H3").Value = "2016-08-22 04:55:07"
$dede.Range("K3").Value = "2016-08-22 04:55:34"
